# Updates cryptos list figures (price / 1h volume change) and reorders two
# coin-name/link pairs, matching the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number (e.g. "0.9987") are
# written with a leading apostrophe so Excel keeps them as text, matching the
# inline-string cells in the workbook (otherwise Excel would silently coerce
# them into numeric cells and drop things like trailing zeros).

$ws.Range("D2").Value = "29.381.67"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "1.841.49"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'239.65"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").Value = "'0.6260"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'0.07444"
$ws.Range("E8").Value = "  -0.73%  "

$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("D10").Value = "'24.77"
$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("D11").Value = "'0.07719"
$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").Value = "1.833.37"
$ws.Range("E12").Value = "  -0.77%  "

$ws.Range("E13").Value = "  -0.61%  "

$ws.Range("D14").Value = "'0.6766"
$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("D15").Value = "'0.00001026"
$ws.Range("E15").Value = "  -2.84%  "

$ws.Range("D16").Value = "'81.73"
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("D17").Value = "'6.239"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").Value = "29.427.02"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("E19").Value = "  +0.97%  "

$ws.Range("E20").Value = "  -0.23%  "

$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").Value = "'7.315"
$ws.Range("E22").Value = "  -2.15%  "

$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").Value = "'158.51"
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("D25").Value = "'8.480"
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("D26").Value = "'0.1353"
$ws.Range("E26").Value = "  -1.60%  "

$ws.Range("D27").Value = "'17.36"
$ws.Range("E27").Value = "  -1.10%  "

$ws.Range("D28").Value = "'0.07331"
$ws.Range("E28").Value = "  +14.74%  "

$ws.Range("D29").Value = "'1.461"
$ws.Range("E29").Value = "  +3.08%  "

$ws.Range("D30").Value = "'1.476"
$ws.Range("E30").Value = "  +0.12%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'4.056"
$ws.Range("E31").Value = "  -1.17%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'4.045"
$ws.Range("E32").Value = "  -1.16%  "

$ws.Range("D33").Value = "'1.819"
$ws.Range("E33").Value = "  -0.30%  "

$ws.Range("D35").Value = "'0.6980"
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("D36").Value = "'2.574"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").Value = "'6.946"
$ws.Range("E37").Value = "  +4.06%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01839"
$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").Value = "'2.817"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("D40").Value = "1.236.08"
$ws.Range("E40").Value = "  -2.33%  "

$ws.Range("D41").Value = "'0.9442"
$ws.Range("E41").Value = "  +3.97%  "

$ws.Range("D42").Value = "'0.9998"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("D43").Value = "1.991.76"
$ws.Range("E43").Value = "  -0.89%  "

$ws.Range("D44").Value = "'101.01"
$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("D45").Value = "'65.64"
$ws.Range("E45").Value = "  -1.02%  "

$ws.Range("E46").Value = "  +2.26%  "

$ws.Range("D47").Value = "'1.721"
$ws.Range("E47").Value = "  -0.83%  "

$ws.Range("D48").Value = "'6.958"
$ws.Range("E48").Value = "  -1.67%  "

$ws.Range("D49").Value = "'8.910"
$ws.Range("E49").Value = "  -1.77%  "

$ws.Range("D51").Value = "'0.3901"
$ws.Range("E51").Value = "  -1.48%  "
